$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B114").Value = 7559468
$ws.Range("E114").Value = 'Liverpool Montevideo'
$ws.Range("F114").Value = 'CA River Plate'
$ws.Range("G114").Value = 2
$ws.Range("H114").Value = 1
$ws.Range("I114").Value = 'H'
$ws.Range("J114").Value = 1.7
$ws.Range("K114").Value = 3
$ws.Range("L114").Value = 5.75
$ws.Range("M114").Value = 1.833
$ws.Range("O114").Value = 4.5
$ws.Range("P114").Value = -0.5
$ws.Range("Q114").Value = 1.925
$ws.Range("R114").Value = 1.925
$ws.Range("S114").Value = 2.25
$ws.Range("T114").Value = 2.025
$ws.Range("U114").Value = 1.825
$ws.Range("V114").Value = 0.833
$ws.Range("W114").Value = -1
$ws.Range("Y114").Value = 0.925
$ws.Range("Z114").Value = -1
$ws.Range("AA114").Value = 1.025
$ws.Range("AB114").Value = -1
$ws.Range("B115").Value = 7559469
$ws.Range("E115").Value = 'Montevideo Wanderers'
$ws.Range("F115").Value = 'Penarol'
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 'D'
$ws.Range("J115").Value = 4.75
$ws.Range("K115").Value = 3.4
$ws.Range("L115").Value = 1.7
$ws.Range("M115").Value = 2.7
$ws.Range("O115").Value = 2.45
$ws.Range("P115").Value = 0
$ws.Range("Q115").Value = 2.05
$ws.Range("R115").Value = 1.8
$ws.Range("S115").Value = 2.5
$ws.Range("T115").Value = 1.975
$ws.Range("U115").Value = 1.875
$ws.Range("V115").Value = -1
$ws.Range("W115").Value = 2.2
$ws.Range("Y115").Value = 0
$ws.Range("Z115").Value = 0
$ws.Range("AA115").Value = -1
$ws.Range("AB115").Value = 0.875
$ws.Range("B117").Value = 7013409
$ws.Range("E117").Value = 'Nacional De Football'
$ws.Range("F117").Value = 'Torque'
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 1
$ws.Range("I117").Value = 'D'
$ws.Range("J117").Value = 1.666
$ws.Range("K117").Value = 3.9
$ws.Range("L117").Value = 4.5
$ws.Range("M117").Value = 1.615
$ws.Range("N117").Value = 4
$ws.Range("O117").Value = 4.75
$ws.Range("P117").Value = -0.75
$ws.Range("Q117").Value = 1.8
$ws.Range("R117").Value = 2.05
$ws.Range("S117").Value = 2.75
$ws.Range("T117").Value = 1.95
$ws.Range("U117").Value = 1.9
$ws.Range("W117").Value = 3
$ws.Range("X117").Value = -1
$ws.Range("Z117").Value = 1.05
$ws.Range("AA117").Value = -1
$ws.Range("AB117").Value = 0.8999999999999999
$ws.Range("B118").Value = 7013885
$ws.Range("E118").Value = 'La Luz'
$ws.Range("F118").Value = 'Atletico Fenix Montevideo'
$ws.Range("J118").Value = 3
$ws.Range("K118").Value = 3
$ws.Range("L118").Value = 2.4
$ws.Range("M118").Value = 2.9
$ws.Range("N118").Value = 2.75
$ws.Range("O118").Value = 2.6
$ws.Range("P118").Value = 0
$ws.Range("Q118").Value = 2.025
$ws.Range("R118").Value = 1.825
$ws.Range("S118").Value = 2
$ws.Range("T118").Value = 2.025
$ws.Range("U118").Value = 1.825
$ws.Range("X118").Value = 1.6
$ws.Range("Z118").Value = 0.825
$ws.Range("AA118").Value = 0
$ws.Range("AB118").Value = 0
$ws.Range("B119").Value = 7013886
$ws.Range("E119").Value = 'Racing Club de Montevideo'
$ws.Range("F119").Value = 'Cerro'
$ws.Range("G119").Value = 0
$ws.Range("I119").Value = 'A'
$ws.Range("J119").Value = 2.25
$ws.Range("K119").Value = 3.1
$ws.Range("L119").Value = 3.25
$ws.Range("M119").Value = 2.25
$ws.Range("N119").Value = 2.875
$ws.Range("O119").Value = 3.5
$ws.Range("P119").Value = -0.25
$ws.Range("Q119").Value = 1.95
$ws.Range("R119").Value = 1.9
$ws.Range("S119").Value = 2
$ws.Range("T119").Value = 1.925
$ws.Range("U119").Value = 1.925
$ws.Range("W119").Value = -1
$ws.Range("X119").Value = 2.5
$ws.Range("Z119").Value = 0.8999999999999999
$ws.Range("AB119").Value = 0.925
$ws.Range("B120").Value = 7013702
$ws.Range("E120").Value = 'Defensor Sporting'
$ws.Range("F120").Value = 'Danubio'
$ws.Range("H120").Value = 2
$ws.Range("J120").Value = 1.8
$ws.Range("K120").Value = 3.6
$ws.Range("L120").Value = 4.2
$ws.Range("M120").Value = 1.8
$ws.Range("N120").Value = 3.6
$ws.Range("O120").Value = 4.2
$ws.Range("P120").Value = -0.75
$ws.Range("Q120").Value = 2.05
$ws.Range("R120").Value = 1.8
$ws.Range("S120").Value = 2.25
$ws.Range("T120").Value = 1.85
$ws.Range("U120").Value = 2
$ws.Range("X120").Value = 3.2
$ws.Range("Z120").Value = 0.8
$ws.Range("AA120").Value = -0.5
$ws.Range("AB120").Value = 0.5
$ws.Range("M210").Value = 1.363
$ws.Range("N210").Value = 4.75
$ws.Range("O210").Value = 7.5
$ws.Range("Q210").Value = 1.875
$ws.Range("R210").Value = 1.975
$ws.Range("T210").Value = 1.9
$ws.Range("U210").Value = 1.95
$ws.Range("M211").Value = 2.75
$ws.Range("N211").Value = 3.1
$ws.Range("O211").Value = 2.625
$ws.Range("S211").Value = 2.25
$ws.Range("M215").Value = 3.1
$ws.Range("N215").Value = 3.2
$ws.Range("O215").Value = 2.1
$ws.Range("P215").Value = 0.25
$ws.Range("Q215").Value = 2.025
$ws.Range("R215").Value = 1.825
$ws.Range("T215").Value = 1.875
$ws.Range("U215").Value = 1.975
